# Generate Report for Handoff
#
# A fresh handback run produced a new working-copy GUID and a new content
# hash for the generated xliff files, so every cell that embeds the old
# identifiers needs to roll over to the new ones, and the two timestamps
# that record when the HO xliff / zh-cn handoff file were produced move
# forward to the new run's generation time.

$wb = $excel.ActiveWorkbook

$oldGuid = "8c9d3c84-5eea-4f92-bc39-03bcf2f58180"
$newGuid = "b3915d0a-ff17-4771-a200-0ec96508a35b"
$oldHash = "c94725b2b865a6b92f73ccf9c3de6db24308d09f"
$newHash = "76ba98f08fdcea797020a99349f42d9554f24a63"

$newMdName = "$newGuid.md"
$newMdDisplay = "e2e\$newGuid.md"

$newZhFile = "$newGuid.$newHash.zh-cn.xlf"
$newDeFile = "$newGuid.$newHash.de-de.xlf"

$newHoDate = "2016-09-01 01:02:05"
$newZhDate = "2016-09-01 01:01:56"

# Replaces the single hyperlink on $cellRef with one that points at the
# same target address but shows $newDisplay as its visible text -
# Hyperlink objects in this object model are immutable in place, so the
# existing link is removed and an equivalent one re-added.
function Update-HyperlinkDisplay($ws, $cellRef, $newDisplay) {
    $addr = $null
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Address
    }
    $ws.Hyperlinks.Delete()
    if ($addr) {
        $ws.Hyperlinks.Add($ws.Range($cellRef), $addr, "", "", $newDisplay) | Out-Null
    }
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value2 = $newMdName
$wsOverview.Range("G2").Value2 = $newHoDate
Update-HyperlinkDisplay $wsOverview "B2" $newMdDisplay

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("G2").Value2 = $newZhFile
$wsZh.Range("H2").Value2 = $newZhDate
Update-HyperlinkDisplay $wsZh "A2" $newMdName

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("G2").Value2 = $newDeFile
$wsDe.Range("H2").Value2 = $newHoDate
Update-HyperlinkDisplay $wsDe "A2" $newMdName
